$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.742.98"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").Value = "3.469.49"
$ws.Range("E3").Value = "  +1.15%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'414.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").Value = "'130.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("D7").Value = "'0.629"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.92%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.727"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("D10").Value = "'0.150"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.94%  "
$ws.Range("D11").Value = "'42.62"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.76%  "
$ws.Range("E12").Value = "  +4.04%  "
$ws.Range("D13").Value = "'0.0000219"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.33%  "
$ws.Range("D14").Value = "4.017.48"
$ws.Range("E14").Value = "  +1.10%  "
$ws.Range("D15").Value = "'0.141"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").Value = "'20.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.74%  "
$ws.Range("D17").Value = "3.468.35"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("D18").Value = "'12.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("E19").Value = "  -1.45%  "
$ws.Range("D20").Value = "62.732.29"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("D21").Value = "'463.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.55%  "
$ws.Range("D22").Value = "'90.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.97%  "
$ws.Range("D23").Value = "'3.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.83%  "
$ws.Range("D24").Value = "'13.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("D25").Value = "'10.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +17.95%  "
$ws.Range("D26").Value = "'3.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.93%  "
$ws.Range("D27").Value = "'33.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.86%  "
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").Value = "'7.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("D30").Value = "'12.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("E31").Value = "  -0.88%  "
$ws.Range("D32").Value = "'0.167"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.61%  "
$ws.Range("E33").Value = "  -1.85%  "
$ws.Range("D34").Value = "'40.76"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.45%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").Value = "'58.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.75%  "
$ws.Range("E37").Value = "  -2.33%  "
$ws.Range("E38").Value = "  +4.60%  "
$ws.Range("D39").Value = "'0.998"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").Value = "'147.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.36%  "
$ws.Range("E41").Value = "  -0.65%  "
$ws.Range("B42").Value = "LidoDAOToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D42").Value = "'3.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.96%  "
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").Value = "'0.321"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.47%  "
$ws.Range("E44").Value = "  +5.12%  "
$ws.Range("D45").Value = "'4.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.31%  "
$ws.Range("D46").Value = "'2.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.68%  "
$ws.Range("D47").Value = "'2.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +14.13%  "
$ws.Range("D48").Value = "0.0₃0556"
$ws.Range("E48").Value = "  +30.31%  "
$ws.Range("D49").Value = "'16.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.50%  "
$ws.Range("D50").Value = "'22.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E51").Value = "  -0.25%  "
